$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 43) down onto the
# two new rows (44-45) so they pick up the same style (borders, fill,
# wrap-text, vertical centering) without introducing new style records.
$ws.Range("A43:G43").Copy()
$ws.Range("A44:G45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new localization rows. Values are written in the same
# interleaved order the original author entered them in (matching the
# resulting shared-string table), rather than strictly left-to-right.
$ws.Range("B44").Value = "tmp game won"
$ws.Range("A44").Value = "Win screen"
$ws.Range("A45").Value = "Death screen"
$ws.Range("C45").Value = "You died!"
$ws.Range("B45").Value = "tmp game lost"
$ws.Range("D44").Value = "Vous avez survécu!"
$ws.Range("C44").Value = "You survived!"
$ws.Range("F45").Value = "死にました！"
$ws.Range("E44").Value = "¡Sobreviviste!"
$ws.Range("F44").Value = "生き残った！"
$ws.Range("G44").Value = "你活下来了！"
$ws.Range("G45").Value = "你死了！"
$ws.Range("E45").Value = "¡Moriste!"
$ws.Range("D45").Value = "Vous avez péri!"

# Match the row height used by the other wrapped two-line rows.
$ws.Rows.Item(44).RowHeight = 28.8
$ws.Rows.Item(45).RowHeight = 28.8

# Reflect the author's final cursor position/selection.
$ws.Range("D43").Select()
